$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = "Oublier la pièce jointe d'un mail, oublier d'enregistrer avant de compiler ou push sur un git, \n oublier de changer de calque avant de modifier un dessin"
$ws.Range("C8").Value = "Utilisation du mauvais raccourci, habitude de disposition des icônes qui mène à des erreurs de clics lors d'un changement, \n fermer trop vite ses notifications"
$ws.Range("B5").Value = "Perturbation cognitive au cours d’une tâche. Vous êtes déconcentré par un événement de l’interface, \n vous avez oublié ce que vous étiez en train de faire."
$ws.Range("C3").Value = "Icône pas claire ou pas assez visible, textes trop longs ou pas assez vulgarisés,  \n îcones ou raccourcis qui sortent des standards"

$ws.Range("C3").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
